$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the monthly time-tracking figures (hours stored as day fractions) ---
# "Navrh" (design) row 3
$ws.Range("G3").Value = 0.36874999999999997
$ws.Range("I3").Value = 0.125

# "Implementace" (implementation) row 4
$ws.Range("G4").Value = 2.4499999999999997
$ws.Range("I4").Value = 1.3784722222222223

# "Psani textu" (writing) row 6
$ws.Range("P6").Value = 1.2916666666666667

# Q column totals (SUM formulas) recalc automatically from the edits above.

# --- View state: zoom + scroll position + selection, matching the author's session ---
$win = $excel.ActiveWindow
$win.Zoom = 250
$win.ScrollRow = 6
$win.ScrollColumn = 1
$ws.Range("L17").Select()

# --- Reposition the time-spent chart lower on the sheet (same size, new anchor) ---
$co = $ws.ChartObjects(1)
$co.Top = 98.05125984251968
$co.Left = 7.944015748031496
